$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A96").Value = "GRT-USD"
